$wb = $excel.ActiveWorkbook

# Rename the second worksheet (tab name truncated to Excel's 31-char limit)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from Primary or Seconda")
$wsInclude.Name = "Include from Classification o"

# Update Metadata sheet values
$wsMeta.Range("B5").Value = "Source Classification Value Set"
$wsMeta.Range("B8").Value = "2021-10-01T15:07:10+00:00"
$wsMeta.Range("B12").Value = "Value set for classifying data origin"

# Update System URI on the renamed "Include from Classification o" sheet
$wsInclude.Range("B4").Value = "http://ibm.com/fhir/cdm/CodeSystem/process-meta-source-classification"
